$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.11353325843811
$ws.Range("B1").Value = 3.076792240142822
$ws.Range("C1").Value = 2.401697158813477
$ws.Range("D1").Value = 2.298105716705322
$ws.Range("E1").Value = 2.230279445648193
